$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.883.93'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.879.56'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '333.06'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4747'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +5.89%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3969'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.72%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '48.15'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08052'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '21.89'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.899.71'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.01%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.967'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.97%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.205'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('E18').Value = '  +1.68%  '
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.29'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.003'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '27.952.86'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.37%  '
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.07'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.73%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.314'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.44%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.121.57'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +3.98%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '157.83'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.11%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.30'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.80%  '
$ws.Range('E29').Value = '  +2.50%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.624'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '122.72'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.47%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9863'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.68%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09579'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.69%  '
$ws.Range('E34').Value = '  -0.49%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.621'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.328'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.56%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06128'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.75%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02264'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.98%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.236'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.85%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.244'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6033'
$ws.Range('D41').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1911'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.38%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '10.38'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.27%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.275'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.40%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5720'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.32'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.414'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.952'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06835'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '113.80'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +5.18%  '
